$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 26
# from serial date 45185 to 45204, preserving existing cell formatting.
for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
